$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores numeric-looking values as plain text in the
# source data (e.g. "42.740.65", "1.01", "72.80"). Force text format before
# assigning so Excel does not auto-convert them to numbers (which would also
# drop meaningful trailing zeros, e.g. "72.80" -> 72.8).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.740.65"
$ws.Range("E2").Value = "  +3.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.267.09"
$ws.Range("E3").Value = "  +3.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.29"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +2.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.51"
$ws.Range("E7").Value = "  +5.93%  "

$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.636"
$ws.Range("E9").Value = "  +7.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.43"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +2.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.05"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.26"
$ws.Range("E13").Value = "  +4.09%  "

$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.612.62"
$ws.Range("E15").Value = "  +3.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.88"
$ws.Range("E16").Value = "  +2.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.877"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.283.83"
$ws.Range("E18").Value = "  +3.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.771.33"
$ws.Range("E19").Value = "  +3.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0991"
$ws.Range("E20").Value = "  +4.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  +1.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.80"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.24"
$ws.Range("E23").Value = "  +9.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.67"
$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.44"
$ws.Range("E28").Value = "  +1.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.65"
$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.10"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.85"
$ws.Range("E32").Value = "  +2.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.26"
$ws.Range("E33").Value = "  +8.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.126"
$ws.Range("E34").Value = "  +4.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0791"
$ws.Range("E35").Value = "  +2.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "30.58"
$ws.Range("E36").Value = "  +18.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("E37").Value = "  +3.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("E38").Value = "  +11.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.70"
$ws.Range("E39").Value = "  +3.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0306"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  +3.90%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.20"
$ws.Range("E42").Value = "  +12.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.97"
$ws.Range("E43").Value = "  +5.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.207"
$ws.Range("E44").Value = "  +7.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.23"
$ws.Range("E45").Value = "  +7.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.63"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.92"
$ws.Range("E47").Value = "  -5.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.38%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.35"
$ws.Range("E51").Value = "  +13.11%  "
